$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two existing July rows with refreshed totals
$ws.Range("B2").Value = 16314.73
$ws.Range("B3").Value = 28535.81

# Insert a new row for July day 3, shifting everything below down by one
$ws.Rows(4).Insert()
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 15955.29
$ws.Range("C4").Value = 7
$ws.Range("D4").Value = 2025
$ws.Range("E4").Value = "07/2025"

# Update the June day-30 total (now shifted down to row 25)
$ws.Range("B25").Value = 114294.26
